# [BI-1613] Update TAF to include term type
#
# Adds a new "Term Type" column (R) to the invalid-fields trait-import
# template: a bold/wrap-text header in R1, and an "Incorrect" value in R3
# (row 2 / the valid-ish record has no term-type error, so it's left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, styled like the other bold headers but on its own
# (wrapped, bold) font so it reads well in a narrower column.
$header = $ws.Range("R1")
$header.Value = "Term Type"
$header.Font.Bold = $true
$header.Font.Size = 11
$header.WrapText = $true

# New invalid-data value for the third data row.
$ws.Range("R3").Value = "Incorrect"

# Reflect where the author ended up looking after adding the column.
$ws.Range("T2").Select() | Out-Null
